$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0.0143196
$ws.Range("F2").Value = 0.008996000000000001
$ws.Range("G2").Value = 0.0233156

$ws.Range("E3").Value = 0.0033008
$ws.Range("F3").Value = 0.029005
$ws.Range("G3").Value = 0.0323058

$ws.Range("E4").Value = 0.0034801
$ws.Range("F4").Value = 0.0189907
$ws.Range("G4").Value = 0.0224708

$ws.Range("E5").Value = 0.0016799
$ws.Range("F5").Value = 0.015511
$ws.Range("G5").Value = 0.0171909

$ws.Range("E6").Value = 0.0017125
$ws.Range("F6").Value = 0.0179787
$ws.Range("G6").Value = 0.0196912

$ws.Range("E7").Value = 0.0021541
$ws.Range("F7").Value = 0.0080588
$ws.Range("G7").Value = 0.0102129

$ws.Range("E8").Value = 0.008999999999999999
$ws.Range("F8").Value = 0.003
$ws.Range("G8").Value = 0.012

$ws.Range("E9").Value = 0.014
$ws.Range("F9").Value = 0.001
$ws.Range("G9").Value = 0.015

$ws.Range("E10").Value = 0.012
$ws.Range("F10").Value = 0.001
$ws.Range("G10").Value = 0.013

$ws.Range("E11").Value = 0.012
$ws.Range("F11").Value = 0.001
$ws.Range("G11").Value = 0.013

$ws.Range("E12").Value = 0.014
$ws.Range("F12").Value = 0.001
$ws.Range("G12").Value = 0.015

$ws.Range("E13").Value = 0.014
$ws.Range("F13").Value = 0.001
$ws.Range("G13").Value = 0.015

$ws.Range("E14").Value = 0.0025464
$ws.Range("F14").Value = 0.0152986
$ws.Range("G14").Value = 0.017845

$ws.Range("E15").Value = 0.0043354
$ws.Range("F15").Value = 0.12559
$ws.Range("G15").Value = 0.1299254

$ws.Range("E16").Value = 0.0030313
$ws.Range("F16").Value = 0.0310622
$ws.Range("G16").Value = 0.0340935

$ws.Range("E17").Value = 0.0034796
$ws.Range("F17").Value = 0.0394506
$ws.Range("G17").Value = 0.0429302

$ws.Range("E18").Value = 0.0281511
$ws.Range("F18").Value = 0.0378699
$ws.Range("G18").Value = 0.066021

$ws.Range("E19").Value = 0.0031338
$ws.Range("F19").Value = 0.0170736
$ws.Range("G19").Value = 0.0202074

$ws.Range("E20").Value = 0.019
$ws.Range("F20").Value = 0.005
$ws.Range("G20").Value = 0.024

$ws.Range("E21").Value = 0.016
$ws.Range("F21").Value = 0.013
$ws.Range("G21").Value = 0.029

$ws.Range("E22").Value = 0.022
$ws.Range("F22").Value = 0.008
$ws.Range("G22").Value = 0.03

$ws.Range("E23").Value = 0.015
$ws.Range("F23").Value = 0.008999999999999999
$ws.Range("G23").Value = 0.024

$ws.Range("E24").Value = 0.021
$ws.Range("F24").Value = 0.011
$ws.Range("G24").Value = 0.032

$ws.Range("E25").Value = 0.019
$ws.Range("F25").Value = 0.015
$ws.Range("G25").Value = 0.034

$ws.Range("E26").Value = 0.0043865
$ws.Range("F26").Value = 0.009944700000000001
$ws.Range("G26").Value = 0.0143312

$ws.Range("E27").Value = 0.002736
$ws.Range("F27").Value = 0.0517939
$ws.Range("G27").Value = 0.0545299

$ws.Range("E28").Value = 0.0024474
$ws.Range("F28").Value = 0.0115114
$ws.Range("G28").Value = 0.0139588

$ws.Range("E29").Value = 0.0024575
$ws.Range("F29").Value = 0.0136257
$ws.Range("G29").Value = 0.0160832

$ws.Range("E30").Value = 0.0032042
$ws.Range("F30").Value = 0.0125908
$ws.Range("G30").Value = 0.015795

$ws.Range("E31").Value = 0.0040657
$ws.Range("F31").Value = 0.0126894
$ws.Range("G31").Value = 0.0167551

$ws.Range("E32").Value = 0.018
$ws.Range("F32").Value = 0.007
$ws.Range("G32").Value = 0.025

$ws.Range("E33").Value = 0.017
$ws.Range("F33").Value = 0.005
$ws.Range("G33").Value = 0.022

$ws.Range("E34").Value = 0.022
$ws.Range("F34").Value = 0.004
$ws.Range("G34").Value = 0.026

$ws.Range("E35").Value = 0.016
$ws.Range("F35").Value = 0.004
$ws.Range("G35").Value = 0.02

$ws.Range("E36").Value = 0.02
$ws.Range("F36").Value = 0.005
$ws.Range("G36").Value = 0.025

$ws.Range("E37").Value = 0.014
$ws.Range("F37").Value = 0.005
$ws.Range("G37").Value = 0.019

$ws.Range("E38").Value = 0.005246
$ws.Range("F38").Value = 0.0003067
$ws.Range("G38").Value = 0.0055527

$ws.Range("E39").Value = 0.004866
$ws.Range("F39").Value = 0.000298
$ws.Range("G39").Value = 0.005163999999999999

$ws.Range("E40").Value = 0.0025165
$ws.Range("F40").Value = 0.000301
$ws.Range("G40").Value = 0.0028175

$ws.Range("E41").Value = 0.0025043
$ws.Range("F41").Value = 0.0003021
$ws.Range("G41").Value = 0.0028064

$ws.Range("E42").Value = 0.0026975
$ws.Range("F42").Value = 0.0004311
$ws.Range("G42").Value = 0.0031286

$ws.Range("E43").Value = 0.0041262
$ws.Range("F43").Value = 0.0002848
$ws.Range("G43").Value = 0.004411

$ws.Range("E44").Value = 0.015
$ws.Range("F44").Value = 0.001
$ws.Range("G44").Value = 0.016

$ws.Range("E45").Value = 0.019
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 0.019

$ws.Range("E46").Value = 0.019
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0.019

$ws.Range("E47").Value = 0.022
$ws.Range("F47").Value = 0.001
$ws.Range("G47").Value = 0.023

$ws.Range("E48").Value = 0.017
$ws.Range("F48").Value = 0.001
$ws.Range("G48").Value = 0.018

$ws.Range("E49").Value = 0.016
$ws.Range("F49").Value = 0.001
$ws.Range("G49").Value = 0.017

